$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.991.03'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.697.46'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.67%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.14%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.124'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.403'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.17%  '

$ws.Range("E12").Value = '  +0.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '30.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.13%  '

$ws.Range("E14").Value = '  +10.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.182.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.823.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.697.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.92%  '

$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '359.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000110'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +15.34%  '

$ws.Range("E26").Value = '  +1.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.173'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.43%  '

$ws.Range("E30").Value = '  +3.77%  '

$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.17%  '

$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '541.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.18%  '

$ws.Range("E33").Value = '  +1.32%  '

$ws.Range("E34").Value = '  +4.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.24%  '

$ws.Range("E36").Value = '  +2.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.60%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.10%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '170.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.89%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0622'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.73'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0267'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.34%  '

$ws.Range("E49").Value = '  +0.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0994'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.50%  '
